# eventbuttons.xlsx edit
# - Adds a new "Artisan Command" entry to the Commands sheet:
#     palette(<int>)  ->  activates palette <int>
#   inserted as a new row right after the "button(<name>)" row (row 87),
#   pushing the existing RC/WebSocket command rows down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")
$ws.Activate()

# Insert a new row at 87 (everything from the old row 87 onward shifts down by one).
$ws.Rows.Item(87).Insert()

# Populate the new row with the palette command documentation.
$ws.Range("B87").Value = "palette(<int>)"
$ws.Range("C87").Value = "activates palette <int>"

# Match the row height used by the neighbouring "Artisan Command" rows.
$ws.Rows.Item(87).RowHeight = 13.8

# Reflect the new selection/scroll position left behind on the Commands sheet.
$ws.Range("B87:C87").Select()
try {
    $excel.ActiveWindow.ScrollRow = 74
    $excel.ActiveWindow.ScrollColumn = 2
} catch {
}

Write-Host "Inserted palette(<int>) command row on Commands sheet."
